# "Updated to 12 tssps on xlsx"
# The pin table used to list 16 Tssp rows (D30:D45) but only the first 12
# (D30:D41) are actually wired to a Tssp sensor - clear the stale "Tssp"
# labels in B43:B46 (D42:D45).
# Also remove the now-unused "Initialization / Pinins / What / Number of
# pins" scratch table that lived in G3:J8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B43:B46").ClearContents()
$ws.Range("G3:J8").ClearContents()

$ws.Range("E39").Select()
